$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pvalue = [double]"7.665750431301586E-14"

# Header row: G1/H1 swap meaning (shared string reorder causes
# column G to now read "significance" and column H to read "pvalue")
$ws.Range("G1").Value = "significance"
$ws.Range("H1").Value = "pvalue"

# Row 2 - group "606"
$ws.Range("A2").Value = "606"
$ws.Range("B2").Value = 9.0
$ws.Range("C2").Value = 29.4
$ws.Range("D2").Value = 3.2667
$ws.Range("E2").Value = 0.1803
$ws.Range("F2").Value = 0.0601
$ws.Range("G2").Value = "c"
$ws.Range("H2").Value = $pvalue

# Row 3 - group "607"
$ws.Range("A3").Value = "607"
$ws.Range("B3").Value = 9.0
$ws.Range("C3").Value = 11.7
$ws.Range("D3").Value = 1.3
$ws.Range("E3").Value = 0.2646
$ws.Range("F3").Value = 0.0882
$ws.Range("G3").Value = "a"
$ws.Range("H3").Value = $pvalue

# Row 4 (new) - group "ZH11"
$ws.Range("A4").Value = "ZH11"
$ws.Range("B4").Value = 9.0
$ws.Range("C4").Value = 22.0
$ws.Range("D4").Value = 2.4444
$ws.Range("E4").Value = 0.3005
$ws.Range("F4").Value = 0.1002
$ws.Range("G4").Value = "b"
$ws.Range("H4").Value = $pvalue
